$wb = $excel.ActiveWorkbook

# --- Fix typo in the Personal Guarantee module name string ---
$wsPG = $wb.Worksheets.Item("KSIDC_Personal_Guarantee")
$wsPG.Range("B2").Value = "Security_Personal_Guarantee_Module"

# --- Add the new "KSIDC_Security_Summary" sheet (after KSIDC_Personal_Guarantee) ---
# Use the Personal Guarantee sheet as a template (same layout / formatting) and
# place the copy immediately after it, i.e. at the end of the workbook.
$wsPG.Copy($null, $wsPG)
$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "KSIDC_Security_Summary"
$wsNew.Range("B2").Value = "Security_Summary_Module"

# --- Update cursor/selection state on the relevant sheets ---
$wsPG.Activate()
$wsPG.Range("B2").Select()

$wsNew.Activate()
$wsNew.Range("B3").Select()
